$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-7: same trawl (14) and species (Notoscopelus kroyerii) as row 2,
# continuing the MINK specimen-code sequence.
$species = "Notoscopelus kroyerii"
$codes = @("MINK003", "MINK004", "MINK005", "MINK006", "MINK007")

$r = 3
foreach ($code in $codes) {
    $ws.Cells.Item($r, 1).Value = 14
    $ws.Cells.Item($r, 2).Value = $species
    $ws.Cells.Item($r, 2).Font.Italic = $true
    $ws.Cells.Item($r, 3).Value = $code
    $ws.Cells.Item($r, 4).Value = "Y"
    $r = $r + 1
}

# Row 8: new trawl (15), new species (Gadiculus thori)
$ws.Cells.Item(8, 1).Value = 15
$ws.Cells.Item(8, 2).Value = "Gadiculus thori"
$ws.Cells.Item(8, 2).Font.Italic = $true
$ws.Cells.Item(8, 3).Value = "MIGT002"
$ws.Cells.Item(8, 4).Value = "Y"

# Page setup, matching a single-page portrait printout on A4/letter-equivalent (paper size 9)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Final selection left on E8, as in the source workbook
[void]$ws.Range("E8").Select()
